# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# This script rebuilds the "Hoja1" account-statement table: it refreshes the
# top summary figures, makes room for the new worker/period rows, fills them
# in with the updated data and leaves the signature footer two blank rows
# below the (now longer) table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Top summary block -----------------------------------------------------
# "VALOR MORA" total amount (row 11)
$ws.Range("E11").Value = 314127

# "Cant. Trabajadores" / "Cant. Periodos" counters (row 13)
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 7

# --- Make room for the new detail rows --------------------------------------
# The table currently has two detail rows (16:17). The refreshed statement
# needs nine (16:24), so insert seven new rows right below the existing ones;
# this pushes the signature footer (previously 22:23) down to 29:30.
$ws.Rows("18:24").Insert()

# Duplicate the "interior" row style (row 16) into row 17 (which used to be
# the closing row) and into all the newly inserted rows, then duplicate the
# "closing" row style (the original row 17, bottom border) into what will
# become the new last detail row (24).
$ws.Range("B17:J17").Copy($ws.Range("B24:J24"))
$ws.Range("B16:J16").Copy($ws.Range("B17:J23"))

# --- Detail rows: worker / period / value data ------------------------------
$rows = @(
    @{ Row=16; Doc="1137222761"; Name="REMBERTO ORTEGA ESPAÑA";     Periodo="2507"; Mora=40000; Salario=1000000 },
    @{ Row=17; Doc="1137222761"; Name="REMBERTO ORTEGA ESPAÑA";     Periodo="2506"; Mora=40000; Salario=1000000 },
    @{ Row=18; Doc="1137222761"; Name="REMBERTO ORTEGA ESPAÑA";     Periodo="2505"; Mora=40000; Salario=1000000 },
    @{ Row=19; Doc="1137222761"; Name="REMBERTO ORTEGA ESPAÑA";     Periodo="2504"; Mora=40000; Salario=1000000 },
    @{ Row=20; Doc="1137222761"; Name="REMBERTO ORTEGA ESPAÑA";     Periodo="2503"; Mora=40000; Salario=1000000 },
    @{ Row=21; Doc="1137222761"; Name="REMBERTO ORTEGA ESPAÑA";     Periodo="2502"; Mora=40000; Salario=1000000 },
    @{ Row=22; Doc="1052984679"; Name="JORGE DAVID PION ALEMAN";    Periodo="1705"; Mora=29509; Salario=737717 },
    @{ Row=23; Doc="1043397732"; Name="YEAN CARLOS CASTILLO CARO";  Periodo="1705"; Mora=29509; Salario=737717 },
    @{ Row=24; Doc="20063315";   Name="CARLOS RAFAEL TORRES PEREZ"; Periodo="1705"; Mora=15109; Salario=781242 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("B$row").Value = "CC"
    $ws.Range("C$row").Value = $r.Doc
    $ws.Range("D$row").Value = $r.Name
    $ws.Range("E$row").Value = $r.Periodo
    $ws.Range("F$row").Value = $r.Mora
    $ws.Range("G$row").Value = $r.Salario
}
